$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "timestamp" column (O) for all data rows (2-73) from
# "2022-07-23 07:02:42" to "2022-07-23 20:59:42".
$ws.Range("O2:O73").Value = "2022-07-23 20:59:42"
